# Weekly price report update:
# Insert a new week's data (2 rows: "Magnum" + "Sin especificar") right after
# the existing row 52, pushing all subsequent weeks down by 2 rows. The final
# two rows that fall off the bottom of the old range simply end up at the new
# bottom of the range (rows 123-124), carrying their original content with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 53:122 down to 55:124 by inserting two blank rows at 53:54.
$ws.Rows("53:54").Insert()

# --- Row 53: Comercializadora del Agro de Limarí, Poroto verde, Magnum ---
$ws.Range("A53").Value = 2
$ws.Range("B53").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C53").Value = "Coquimbo"
$ws.Range("D53").Value2 = 44546
$ws.Range("D53").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E53").Value = 4
$ws.Range("F53").Value = 100112031
$ws.Range("G53").Value = "Poroto verde"
$ws.Range("H53").Value = "Magnum"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 600
$ws.Range("K53").Value = 13000
$ws.Range("L53").Value = 14000
$ws.Range("M53").Value = 13500
$ws.Range("N53").Value = "$/malla 25 kilos"
$ws.Range("O53").Value = "Provincia de Limarí"
$ws.Range("P53").Value = 540
$ws.Range("Q53").Value = 25
$ws.Range("R53").Value = "Hortaliza"

# --- Row 54: Comercializadora del Agro de Limarí, Poroto verde, Sin especificar ---
$ws.Range("A54").Value = 2
$ws.Range("B54").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C54").Value = "Coquimbo"
$ws.Range("D54").Value2 = 44546
$ws.Range("D54").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E54").Value = 4
$ws.Range("F54").Value = 100112031
$ws.Range("G54").Value = "Poroto verde"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 400
$ws.Range("K54").Value = 22000
$ws.Range("L54").Value = 24000
$ws.Range("M54").Value = 23000
$ws.Range("N54").Value = "$/malla 25 kilos"
$ws.Range("O54").Value = "Provincia de Limarí"
$ws.Range("P54").Value = 920
$ws.Range("Q54").Value = 25
$ws.Range("R54").Value = "Hortaliza"
